# Preliminary work so that readInputSheet reads ALL files in the test_files
# directory, not just the sixteen_tests, etc.
#
# optimization_parameters sheet: "Model" becomes "production_function", and a
# new "L_curve" parameter row is inserted right after it. The old "Deletion"
# row (Strain sub-table) is removed. The previously-active tab
# (optimization_parameters) is no longer the active tab; wt_log2_expression
# becomes the active tab instead, and the selection on
# optimization_parameters moves to E9.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("optimization_parameters")

# Drop the stray duplicated "value" header cells C1:F1 (row 1 only really
# needs A1/B1 -- "optimization_parameter" / "value")
$wsParams.Range("C1:F1").ClearContents()

# "Model" -> "production_function" (value/style of B8 "Sigmoid" unchanged)
$wsParams.Range("A8").Value = "production_function"

# Insert a new row for "L_curve" right after "production_function"
[void]$wsParams.Rows.Item(9).Insert()
$wsParams.Range("A9").Value = "L_curve"
$wsParams.Range("B9").Value = 0
$wsParams.Range("B9").NumberFormat = "0.00E+00"

# Remove the old "Deletion" row (now shifted down to row 17 after the insert)
[void]$wsParams.Rows.Item(17).Delete()

# optimization_parameters keeps its own cursor position, now at E9 (select
# this first -- selecting a range on a sheet also activates that sheet, so
# the final Activate() below must come after this to end up on the right tab)
[void]$wsParams.Range("E9").Select()

# Move the active tab from optimization_parameters to wt_log2_expression
$wsExpr = $wb.Worksheets.Item("wt_log2_expression")
[void]$wsExpr.Activate()
